# Append: 2025-10-19 06:25 JST
# The listing sheet ("ランサーズ", the first/active sheet) is refreshed with a
# new scrape: the previous 13 data rows (rows 2-14) are replaced by 5 new
# data rows (rows 2-6), column widths for B and H are widened, and the
# hyperlinks on column F are rebuilt to point at the new listing URLs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Hyperlinks on this sheet are rebuilt from scratch -----------------
# (Deleting via any Range's .Hyperlinks collection clears every hyperlink on
# the sheet in this engine, so do it once up-front before re-adding them.)
$ws.Hyperlinks.Delete()

# --- 2. Drop the old rows 7-14 (only 5 data rows remain after the refresh) -
$ws.Rows("7:14").Delete()

# --- 3. Row 2 -------------------------------------------------------------
$ws.Range("A2").Value = "2025-10-19 06:25:53"
$ws.Range("B2").Value = "GoogleスプレッドシートとMetaAPIを利用したFXトレード大会ランキングの自動化システム開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5416128"
$ws.Range("G2").Value = 363
$ws.Range("H2").Value = "🔥API ◆開発,システム開発"

# --- 4. Row 3 -------------------------------------------------------------
$ws.Range("A3").Value = "2025-10-19 06:25:53"
$ws.Range("B3").Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Range("G3").Value = 135
$ws.Range("H3").Value = "◆ツール,スクレイピング ◇サイト"

# --- 5. Row 4 -------------------------------------------------------------
$ws.Range("A4").Value = "2025-10-19 06:25:53"
$ws.Range("B4").Value = "イベント出店者管理用ウェブアプリ開発依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5416005"
$ws.Range("G4").Value = 123
$ws.Range("H4").Value = "◆開発 ◇アプリ"

# --- 6. Row 5 -------------------------------------------------------------
$ws.Range("A5").Value = "2025-10-19 06:25:53"
$ws.Range("B5").Value = "【恋愛診断】フルスクラッチ開発・運用サポート募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5415986"
$ws.Range("G5").Value = 75
$ws.Range("H5").Value = "◆開発"

# --- 7. Row 6 (no skill-summary cell, same as before the refresh) ---------
$ws.Range("A6").Value = "2025-10-19 06:25:53"
$ws.Range("B6").Value = "【高額成功報酬】レガシー基幹システムのバイナリ解析とパッチ作成"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5415980"
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = ""

# --- 8. Re-create the hyperlinks for the URL column (F2:F6) ---------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5416128")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5251319")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5416005")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5415986")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5415980")

# --- 9. Widen columns B and H -----------------------------------------
# Excel's ColumnWidth setter pads by ~5px (0.8333 chars at the default
# font), so back that padding out to land exactly on the target stored
# widths of 52 and 19 characters.
$ws.Columns.Item(2).ColumnWidth = 51.16666666666666
$ws.Columns.Item(8).ColumnWidth = 18.16666666666667
